# Revert "Merging 0.1.8 w VitalSigns"
#
# Restores the Metadata sheet to its pre-merge values (version, status,
# date, contact info), removes the extra "Jurisdiction" row that the
# merge introduced, and renames the two "Include ValueSet #N" sheets
# back to "Include ValueSets" / "Include ValueSets 2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update metadata values on the "Metadata" sheet ---
$ws.Range("B3").Value  = "0.1.6"
$ws.Range("B6").Value  = "active"
$ws.Range("B8").Value  = "2023-05-05T10:50:04-05:00"
$ws.Range("B10").Value = "No display for ContactDetail"
$ws.Range("B11").Value = "No display for ContactDetail"

# --- Remove the "Jurisdiction" row (old row 12); rows below shift up ---
$ws.Range("A12:B12").EntireRow.Delete()

# --- Rename the "Include ValueSet #N" sheets ---
$wb.Worksheets.Item("Include ValueSet #0").Name = "Include ValueSets"
$wb.Worksheets.Item("Include ValueSet #1").Name = "Include ValueSets 2"
